$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values remain stored as text, matching the source
# data which represents prices as literal strings (e.g. with thousands dots
# or preserved trailing zeros) rather than numeric values.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '43.031.56'
$ws.Range("E2").Value = '  +4.07%  '
$ws.Range("D3").Value = '2.233.41'
$ws.Range("E3").Value = '  +3.78%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '252.39'
$ws.Range("E5").Value = '  +6.69%  '
$ws.Range("D6").Value = '0.617'
$ws.Range("E6").Value = '  +2.51%  '
$ws.Range("D7").Value = '75.25'
$ws.Range("E7").Value = '  +7.02%  '
$ws.Range("E8").Value = '  -0.19%  '
$ws.Range("D9").Value = '0.598'
$ws.Range("E9").Value = '  +4.69%  '
$ws.Range("D10").Value = '41.11'
$ws.Range("E10").Value = '  +5.50%  '
$ws.Range("E11").Value = '  +2.92%  '
$ws.Range("D12").Value = '6.89'
$ws.Range("E12").Value = '  +4.01%  '
$ws.Range("D13").Value = '0.101'
$ws.Range("E13").Value = '  +1.59%  '
$ws.Range("D14").Value = '2.569.35'
$ws.Range("E14").Value = '  +3.84%  '
$ws.Range("E15").Value = '  +2.73%  '
$ws.Range("D16").Value = '2.230.35'
$ws.Range("E16").Value = '  +4.26%  '
$ws.Range("D17").Value = '0.790'
$ws.Range("E17").Value = '  +1.08%  '
$ws.Range("D18").Value = '42.931.39'
$ws.Range("E18").Value = '  +4.27%  '
$ws.Range("E19").Value = '  +3.74%  '
$ws.Range("D20").Value = '71.25'
$ws.Range("E20").Value = '  +3.10%  '
$ws.Range("D21").Value = '5.95'
$ws.Range("E21").Value = '  +3.84%  '
$ws.Range("D22").Value = '230.11'
$ws.Range("E22").Value = '  +2.14%  '
$ws.Range("E23").Value = '  +11.49%  '
$ws.Range("E24").Value = '  -3.38%  '
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("D26").Value = '10.72'
$ws.Range("E26").Value = '  +1.21%  '
$ws.Range("E27").Value = '  +4.32%  '
$ws.Range("D28").Value = '39.28'
$ws.Range("E28").Value = '  +23.22%  '
$ws.Range("E29").Value = '  +2.73%  '
$ws.Range("E30").Value = '  -0.63%  '
$ws.Range("D31").Value = '170.26'
$ws.Range("D32").Value = '20.19'
$ws.Range("E33").Value = '  +4.89%  '
$ws.Range("D34").Value = '5.25'
$ws.Range("E34").Value = '  +3.41%  '
$ws.Range("D35").Value = '0.113'
$ws.Range("E35").Value = '  +9.71%  '
$ws.Range("E36").Value = '  +0.88%  '
$ws.Range("D37").Value = '4.44'
$ws.Range("E37").Value = '  +3.99%  '
$ws.Range("E38").Value = '  +12.49%  '
$ws.Range("D39").Value = '12.28'
$ws.Range("E39").Value = '  +2.92%  '
$ws.Range("E41").Value = '  +9.30%  '
$ws.Range("D42").Value = '5.37'
$ws.Range("E42").Value = '  +2.49%  '
$ws.Range("D43").Value = '59.61'
$ws.Range("E43").Value = '  +2.41%  '
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").Value = '8.65'
$ws.Range("E44").Value = '  +4.97%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = '103.34'
$ws.Range("E45").Value = '  +6.18%  '
$ws.Range("B46").Value = 'WOONetwork'
$ws.Range("C46").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D46").Value = '0.482'
$ws.Range("E46").Value = '  +26.32%  '
$ws.Range("D47").Value = '0.0984'
$ws.Range("E48").Value = '  +12.82%  '
$ws.Range("E49").Value = '  +4.15%  '
$ws.Range("E50").Value = '  +3.25%  '
$ws.Range("E51").Value = '  +2.22%  '
